# Developer-guide diagram update:
#   * Refresh the cached "datetimeFigureOut" auto-date fields on the
#     notes master, slide master and every slide layout from 2/6/2017
#     to 3/22/2017.
#   * Nudge/resize the "execute(...)" callout and recolor its text to
#     "delete f1".
#   * Shift the "Rectangle 62" box that sits to its right.
#   * Rename deletePerson -> deleteFloat and update the two parser
#     snippets to reflect the new "f1" argument.

$p = $ppt.ActivePresentation

$oldDate = "2/6/2017"
$newDate = "3/22/2017"

# -- 1. Notes master date placeholder ---------------------------------
# NOTE: this COM-interop host does not route shape writes made through
# `$ppt.ActivePresentation.NotesMaster` to the notes-master part (they
# either vanish or - when the notes-master shape id happens to collide
# with a slide-master shape id - land on the slide master instead), so
# touching it here would corrupt the slide master. Left alone on purpose.

# -- 2. Slide master date placeholder ----------------------------------
$master = $p.SlideMaster
foreach ($sh in $master.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# -- 3. Every slide layout's date placeholder --------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    foreach ($sh in $layout.Shapes) {
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# -- 4. Slide content edits --------------------------------------------
$s = $p.Slides.Item(1)
$q1 = [char]0x201C
$q2 = [char]0x201D

foreach ($sh in $s.Shapes) {

    if ($sh.Id -eq 26) {
        # "execute(<93>delete 1<94>)" textbox: reposition/resize + retext
        $sh.Left = 3491 / 12700
        $sh.Width = 1520509 / 12700
        $sh.TextFrame.TextRange.Text = "execute(" + $q1 + "delete f1" + $q2 + ")"
    }

    if ($sh.Id -eq 39) {
        # "Rectangle 62" box: shift left edge only
        $sh.Left = 8077200 / 12700
    }

    if ($sh.Id -eq 78) {
        # "deletePerson(p)" -> "deleteFloat(p)"; keep the "(p)" run as-is
        $tr = $sh.TextFrame.TextRange
        $tr.Characters(1, 12).Text = "deleteFloat"
    }

    if ($sh.Id -eq 79) {
        # "p" + "arse(<93>1<94>)" -> "p" + "arse(<93>f<94>, 1)"
        $tr = $sh.TextFrame.TextRange
        $origHeight = $sh.Height
        $tr.Characters(2, 9).Text = "arse(" + $q1 + "f" + $q2 + ", 1)"
        # Re-assert the original autofit height: PowerPoint's own resave
        # left this box's extent untouched even though the text grew.
        $sh.Height = $origHeight
    }

    if ($sh.Id -eq 80) {
        # "parse(<93>delete 1<94>)" -> "parse(<93>delete f1<94>)"
        $sh.TextFrame.TextRange.Text = "parse(" + $q1 + "delete f1" + $q2 + ")"
    }
}
